$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Column A (IES) values for the rows touched by this edit
$ies1 = "UNIVERSIDADES"
$ies2 = "INSTITUCION UNIVERSITARIA"

# Existing rows 2-6 shift to new ids/universities and pick up a new
# "UNIVERSIDADES" label in column A (was PONTIFICIA UNIVERSIDAD JAVERIANA's
# own shared string before, now re-pointed at the UNIVERSIDADES entry).
$data = @(
    @(2,  $ies1, 12385, "UNIVERSIDAD DE LOS ANDES ",            "PRIVADA", "UNIVERSIDAD", "SI", 1),
    @(3,  $ies1, 12386, "UNIVERSIDAD NACIONAL DE COLOMBIA",     "OFICIAL", "UNIVERSIDAD", "SI", 2),
    @(4,  $ies1, 12387, "UNIVERSIDAD DEL TOLIMA",                "PRIVADA", "UNIVERSIDAD", "SI", 3),
    @(5,  $ies1, 12388, "UNIVERSIDAD DEL NORTE",                 "PRIVADA", "UNIVERSIDAD", "SI", 4),
    @(6,  $ies1, 12389, "PONTIFICIA UNIVERSIDAD JAVERIANA",      "PRIVADA", "UNIVERSIDAD", "SI", 5),
    @(7,  $ies1, 12390, "UNIVERSIDAD GRAN COLOMBIANA",           "PRIVADA", "UNIVERSIDAD", "SI", 6),
    @(8,  $ies1, 12391, "UNIVERSIDAD PARA TODOS",                "OFICIAL", "UNIVERSIDAD", "NO", 7),
    @(9,  $ies1, 12392, "UNIVERSIDAD DE TUNJA",                  "PRIVADA", "UNIVERSIDAD", "SI", 8),
    @(10, $ies1, 12393, "UNIVERSIDAD DEL IBAGE",                 "PRIVADA", "UNIVERSIDAD", "SI", 9),
    @(11, $ies1, 12394, "PONTIFICIA UNIVERSIDAD DEL VALLE",      "PRIVADA", "UNIVERSIDAD", "NO", 10),
    @(12, $ies1, 12395, "UNIVERSIDAD DE LOS FLORES",             "PRIVADA", "UNIVERSIDAD", "SI", 11),
    @(13, $ies2, 12396, "UNIVERSIDAD DE CALI",                   "OFICIAL", "UNIVERSIDAD", "SI", 12),
    @(14, $ies2, 12397, "UNIVERSIDAD DE MEDELLIN",                "PRIVADA", "UNIVERSIDAD", "SI", 13),
    @(15, $ies2, 12398, "UNIVERSIDAD DEL NORTE DEL VALLE",        "PRIVADA", "UNIVERSIDAD", "SI", 14),
    @(16, $ies2, 12399, "PONTIFICIA UNIVERSIDAD PACIFICO",        "PRIVADA", "UNIVERSIDAD", "SI", 15),
    @(17, $ies2, 12400, "UNIVERSIDAD DEL PACIFICO NORTE",         "PRIVADA", "UNIVERSIDAD", "SI", 16),
    @(18, $ies2, 12401, "UNIVERSIDAD PLAN NACIONAL",              "OFICIAL", "UNIVERSIDAD", "SI", 17),
    @(19, $ies2, 12402, "UNIVERSIDAD DE LAS MARIAS",              "PRIVADA", "UNIVERSIDAD", "SI", 18),
    @(20, $ies2, 12403, "UNIVERSIDAD DEL PACICO",                 "PRIVADA", "UNIVERSIDAD", "SI", 19),
    @(21, $ies2, 12404, "PONTIFICIA UNIVERSIDAD DE LA REGION",    "PRIVADA", "UNIVERSIDAD", "SI", 20)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $data[$i]
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 1).HorizontalAlignment = -4131
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 2).Font.Size = 12
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]
    $ws.Cells.Item($row, 7).Value = $r[7]
}

$ws.Columns.Item(1).ColumnWidth = 26.85546875
$ws.Range("A13:A21").Select()
